$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to match the data it holds.
$ws.Name = "Liverpool Stats"

# Turn the A1:M2 range into a proper Excel Table ("Table1") so the stats
# can be filtered/sorted like a real data set.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:M2"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Widen several columns that were previously only "best fit" so the
# (now table-ified) headers/values are comfortably readable.
$ws.Columns(2).ColumnWidth = 14.54296875
$ws.Columns(4).ColumnWidth = 15.90625
$ws.Columns(8).ColumnWidth = 10.54296875
$ws.Columns(9).ColumnWidth = 14.08984375
$ws.Columns(10).ColumnWidth = 15.81640625
$ws.Columns(11).ColumnWidth = 11.54296875
$ws.Columns(12).ColumnWidth = 11.36328125
$ws.Columns(13).ColumnWidth = 13.6328125

# Scroll back to the default top-left and move the selection.
[void]$ws.Range("D19").Select()
